$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 107
$ws.Range("D107").Value = 0.5
$ws.Range("E107").Value = 5.9
$ws.Range("G107").Value = 92.59999999999999
$ws.Range("H107").Value = 0.01
$ws.Range("I107").Value = 0.08
$ws.Range("J107").Value = 0.01
$ws.Range("K107").Value = 0.02
$ws.Range("L107").Value = 0.01

# Row 108
$ws.Range("D108").Value = 1
$ws.Range("E108").Value = 8.5
$ws.Range("G108").Value = 89.5
$ws.Range("H108").Value = 0.01
$ws.Range("I108").Value = 0.11
$ws.Range("J108").Value = 0.02
$ws.Range("K108").Value = 0.04
$ws.Range("L108").Value = 0.01

# Row 109
$ws.Range("D109").Value = 0.2
$ws.Range("E109").Value = 3.2
$ws.Range("G109").Value = 95.59999999999999
$ws.Range("I109").Value = 0.06
$ws.Range("K109").Value = 0.01
